$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Drop the stray "_GoBack" bookmark that used to sit right after
#    "(requires MS Office 2016)" in the Mac-accessibility-checker
#    paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Turn the trailing empty paragraph into the new "Screenreaders Info"
#    Heading 2, and append the four new reference paragraphs after it
#    (raw-XML splice keeps exact run/paragraph formatting + wires up the
#    three new hyperlink relationships in one shot).
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Screenreaders</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Info</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:lastRenderedPageBreak/><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "http://webaccess.berkeley.edu/file/screen-reader-demo" </w:instrText></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:t>“Screen Reader Demo,” Berkeley Web Access @ UC Berkeley https://webaccess.berkeley.edu/file/screen-reader-demo</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:hyperlink r:id="rId13" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve">"Screen Reader Simulation,” </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:t>WebAIM</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:br/><w:t>http://webaim.org/simulations/screenreader</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:hyperlink r:id="rId14" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:t>“Screen reader demos, research and resources,” by Alistair Duggan http://alistairduggin.co.uk/blog/screenreader-resources/</w:t></w:r></w:hyperlink></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:hyperlink r:id="rId15" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve">“Designing for Screen Reader Compatibility,” </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:t>WebAIM</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> http://webaim.org/techniques/screenreader/</w:t></w:r></w:hyperlink></w:p></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId13" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://webaim.org/simulations/screenreader" TargetMode="External"/><Relationship Id="rId14" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://alistairduggin.co.uk/blog/screenreader-resources/" TargetMode="External"/><Relationship Id="rId15" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://webaim.org/techniques/screenreader/" TargetMode="External"/></Relationships></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($fragment)

# ---------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark around the four reference
#    paragraphs (everything after the new heading).
# ---------------------------------------------------------------------
$linksStart = $d.Paragraphs.Item(13).Range.Start
$linksEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$bmRange = $d.Range($linksStart, $linksEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
